# amazonSteps.xlsx update:
# - new changes for multiple retailers and allure reporting
# Rework Sheet1's test-step table (fewer, different rows), add a blank
# "Sheet2", widen column B, bump E2 to 5000, and give the "empty separator"
# rows in column A a bold style (style index 3).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Remove old hyperlinks before wiping cell contents ---------------------
$ws1.Hyperlinks.Delete()

# --- Wipe all existing cell content/formatting on Sheet1 --------------------
$ws1.Cells.Clear()

# --- Column B is much wider now ---------------------------------------------
$ws1.Columns("B").ColumnWidth = 73.109375

# Helper alignment constants (xlHAlignLeft = -4131, xlVAlignTop = -4160)
$hLeft = -4131
$vTop = -4160

# --- Row 1: header -----------------------------------------------------------
$ws1.Range("A1").Value = "action"
$ws1.Range("B1").Value = "locator"
$ws1.Range("C1").Value = "value"
$ws1.Range("D1").Value = "waitBefore"
$ws1.Range("E1").Value = "waitAfter"

# --- Row 2 --------------------------------------------------------------------
$ws1.Range("A2").Value = "goto"
$ws1.Range("B2").Value = "https://www.amazon.com/"
$ws1.Range("D2").Value = 1000
$ws1.Range("E2").Value = 5000

# --- Row 3 (bold "waitfortext" marker in column A) -----------------------------
$ws1.Range("A3").Value = "waitfortext"
$ws1.Range("A3").Font.Bold = $true
$ws1.Range("B3").Value = "Hello, Sign in"

# --- Row 4 (write "Automation" before the span locator so the shared-string
#            table appends them in the same order the real edit produced) ------
$ws1.Range("A4").Value = "assert"
$ws1.Range("C4").Value = "Automation"
$ws1.Range("B4").Value = "span#nav-link-accountList-nav-line-1"
$ws1.Range("D4").Value = 1000

# --- Rows 5-9: bold, empty placeholder cells in column A ------------------------
$ws1.Range("A5").Font.Bold = $true
$ws1.Range("A6").Font.Bold = $true
$ws1.Range("A7").Font.Bold = $true
$ws1.Range("A8").Font.Bold = $true
$ws1.Range("A9").Font.Bold = $true

# --- Row 10: empty cell carrying the "Hyperlink" style (no hyperlink) ----------
$ws1.Range("B10").Style = "Hyperlink"

$ws1.Range("A1").Select()

# --- Add Sheet2 right after Sheet1, blank, selection on B3 ----------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Range("B3").Select()

$ws1.Activate()
